$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) contains values such as "0.5227" or "1.003" that
# Excel would otherwise auto-convert to numbers when assigned through
# .Value. Force the whole column to Text format first so every new value
# is stored as a string, matching the source data feed.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.122.28'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '1.668.68'
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").Value = '210.91'
$ws.Range("E5").Value = '  -2.31%  '
$ws.Range("D6").Value = '0.5227'
$ws.Range("E6").Value = '  -0.86%  '
$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  -0.30%  '
$ws.Range("D8").Value = '0.2624'
$ws.Range("E8").Value = '  -2.33%  '
$ws.Range("D9").Value = '0.06334'
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("D10").Value = '21.22'
$ws.Range("E10").Value = '  -1.08%  '
$ws.Range("D11").Value = '0.07553'
$ws.Range("E11").Value = '  -0.96%  '
$ws.Range("D12").Value = '1.677.90'
$ws.Range("E12").Value = '  -0.41%  '
$ws.Range("D13").Value = '4.433'
$ws.Range("E13").Value = '  -2.04%  '
$ws.Range("D14").Value = '0.5486'
$ws.Range("E14").Value = '  -4.53%  '
$ws.Range("D15").Value = '0.000008049'
$ws.Range("E15").Value = '  -2.31%  '
$ws.Range("D16").Value = '66.47'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '26.176.47'
$ws.Range("E17").Value = '  -0.20%  '
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("D19").Value = '4.753'
$ws.Range("E19").Value = '  -2.36%  '
$ws.Range("D20").Value = '187.95'
$ws.Range("E20").Value = '  -0.93%  '
$ws.Range("D21").Value = '10.30'
$ws.Range("E21").Value = '  -4.00%  '
$ws.Range("D22").Value = '6.237'
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("E23").Value = '  -0.30%  '
$ws.Range("D24").Value = '149.53'
$ws.Range("E24").Value = '  +0.26%  '
$ws.Range("D25").Value = '0.1241'
$ws.Range("E25").Value = '  -1.46%  '
$ws.Range("D26").Value = '7.479'
$ws.Range("E26").Value = '  -3.04%  '
$ws.Range("D27").Value = '15.81'
$ws.Range("E27").Value = '  -0.19%  '
$ws.Range("D28").Value = '0.06334'
$ws.Range("E28").Value = '  -1.07%  '
$ws.Range("D29").Value = '1.355'
$ws.Range("E29").Value = '  -1.54%  '
$ws.Range("D30").Value = '1.283'
$ws.Range("E30").Value = '  -2.38%  '
$ws.Range("D31").Value = '3.528'
$ws.Range("E31").Value = '  -1.00%  '
$ws.Range("D32").Value = '3.416'
$ws.Range("E32").Value = '  -4.25%  '
$ws.Range("D33").Value = '1.648'
$ws.Range("E33").Value = '  -1.96%  '
$ws.Range("D34").Value = '1.005'
$ws.Range("E34").Value = '  -1.28%  '
$ws.Range("D35").Value = '0.6025'
$ws.Range("E35").Value = '  -1.46%  '
$ws.Range("B36").Value = 'MXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D36").Value = '2.763'
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '2.396'
$ws.Range("E37").Value = '  -1.06%  '
$ws.Range("D38").Value = '1.118.04'
$ws.Range("E38").Value = '  +2.13%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01613'
$ws.Range("E39").Value = '  -1.69%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '6.066'
$ws.Range("E40").Value = '  -1.70%  '
$ws.Range("D41").Value = '0.8643'
$ws.Range("E41").Value = '  -1.95%  '
$ws.Range("D42").Value = '1.003'
$ws.Range("E42").Value = '  -0.64%  '
$ws.Range("D43").Value = '100.46'
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("D44").Value = '1.822.02'
$ws.Range("E44").Value = '  -0.64%  '
$ws.Range("D45").Value = '0.00000000108'
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").Value = '55.45'
$ws.Range("E46").Value = '  -3.54%  '
$ws.Range("D47").Value = '1.002'
$ws.Range("E47").Value = '  +0.21%  '
$ws.Range("D48").Value = '8.063'
$ws.Range("E48").Value = '  -0.30%  '
$ws.Range("D49").Value = '0.05237'
$ws.Range("E49").Value = '  -0.57%  '
$ws.Range("D50").Value = '0.4239'
$ws.Range("E50").Value = '  -0.96%  '
$ws.Range("D51").Value = '5.921'
$ws.Range("E51").Value = '  -1.56%  '

# Restore the default (Normal) style on the Price column so no stray
# number-format metadata remains attached to the cells.
$priceRange.Style = "Normal"
